$d = $word.ActiveDocument

function Replace-Once($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# Phase 1: replace each original value with a unique, collision-free placeholder
# so that later phase-1 replacements cannot accidentally match a value produced
# by an earlier phase-1 replacement (the Find engine always scans from the start
# of the document for the first match).
Replace-Once $d.Content "2025-08-23 Saturday" "@@PLACEHOLDER_0@@"
Replace-Once $d.Content "52÷3=17, 1" "@@PLACEHOLDER_1@@"
Replace-Once $d.Content "80÷2=40, 0" "@@PLACEHOLDER_2@@"
Replace-Once $d.Content "78÷7=11, 1" "@@PLACEHOLDER_3@@"
Replace-Once $d.Content "74÷2=37, 0" "@@PLACEHOLDER_4@@"
Replace-Once $d.Content "54÷6=9, 0" "@@PLACEHOLDER_5@@"
Replace-Once $d.Content "72÷9=8, 0" "@@PLACEHOLDER_6@@"
Replace-Once $d.Content "75÷4=18, 3" "@@PLACEHOLDER_7@@"
Replace-Once $d.Content "64÷8=8, 0" "@@PLACEHOLDER_8@@"
Replace-Once $d.Content "82÷6=13, 4" "@@PLACEHOLDER_9@@"
Replace-Once $d.Content "36÷4=9, 0" "@@PLACEHOLDER_10@@"
Replace-Once $d.Content "17÷5=3, 2" "@@PLACEHOLDER_11@@"
Replace-Once $d.Content "44÷7=6, 2" "@@PLACEHOLDER_12@@"
Replace-Once $d.Content "48÷9=5, 3" "@@PLACEHOLDER_13@@"
Replace-Once $d.Content "28÷7=4, 0" "@@PLACEHOLDER_14@@"
Replace-Once $d.Content "99÷7=14, 1" "@@PLACEHOLDER_15@@"
Replace-Once $d.Content "85÷4=21, 1" "@@PLACEHOLDER_16@@"
Replace-Once $d.Content "20÷5=4, 0" "@@PLACEHOLDER_17@@"
Replace-Once $d.Content "11÷7=1, 4" "@@PLACEHOLDER_18@@"
Replace-Once $d.Content "74÷7=10, 4" "@@PLACEHOLDER_19@@"
Replace-Once $d.Content "54÷8=6, 6" "@@PLACEHOLDER_20@@"
Replace-Once $d.Content "62÷9=6, 8" "@@PLACEHOLDER_21@@"
Replace-Once $d.Content "15÷6=2, 3" "@@PLACEHOLDER_22@@"
Replace-Once $d.Content "58÷8=7, 2" "@@PLACEHOLDER_23@@"
Replace-Once $d.Content "32÷9=3, 5" "@@PLACEHOLDER_24@@"
Replace-Once $d.Content "26÷3=8, 2" "@@PLACEHOLDER_25@@"

# Phase 2: replace each placeholder with its final value. Placeholders are
# unique strings that cannot collide with each other or with any real content.
Replace-Once $d.Content "@@PLACEHOLDER_0@@" "2025-08-24 Sunday"
Replace-Once $d.Content "@@PLACEHOLDER_1@@" "10÷8=1, 2"
Replace-Once $d.Content "@@PLACEHOLDER_2@@" "63÷6=10, 3"
Replace-Once $d.Content "@@PLACEHOLDER_3@@" "20÷5=4, 0"
Replace-Once $d.Content "@@PLACEHOLDER_4@@" "96÷6=16, 0"
Replace-Once $d.Content "@@PLACEHOLDER_5@@" "35÷4=8, 3"
Replace-Once $d.Content "@@PLACEHOLDER_6@@" "93÷5=18, 3"
Replace-Once $d.Content "@@PLACEHOLDER_7@@" "30÷5=6, 0"
Replace-Once $d.Content "@@PLACEHOLDER_8@@" "80÷2=40, 0"
Replace-Once $d.Content "@@PLACEHOLDER_9@@" "76÷8=9, 4"
Replace-Once $d.Content "@@PLACEHOLDER_10@@" "43÷8=5, 3"
Replace-Once $d.Content "@@PLACEHOLDER_11@@" "89÷9=9, 8"
Replace-Once $d.Content "@@PLACEHOLDER_12@@" "95÷8=11, 7"
Replace-Once $d.Content "@@PLACEHOLDER_13@@" "47÷2=23, 1"
Replace-Once $d.Content "@@PLACEHOLDER_14@@" "38÷9=4, 2"
Replace-Once $d.Content "@@PLACEHOLDER_15@@" "76÷6=12, 4"
Replace-Once $d.Content "@@PLACEHOLDER_16@@" "77÷5=15, 2"
Replace-Once $d.Content "@@PLACEHOLDER_17@@" "26÷4=6, 2"
Replace-Once $d.Content "@@PLACEHOLDER_18@@" "29÷6=4, 5"
Replace-Once $d.Content "@@PLACEHOLDER_19@@" "84÷3=28, 0"
Replace-Once $d.Content "@@PLACEHOLDER_20@@" "97÷3=32, 1"
Replace-Once $d.Content "@@PLACEHOLDER_21@@" "24÷8=3, 0"
Replace-Once $d.Content "@@PLACEHOLDER_22@@" "85÷7=12, 1"
Replace-Once $d.Content "@@PLACEHOLDER_23@@" "85÷3=28, 1"
Replace-Once $d.Content "@@PLACEHOLDER_24@@" "26÷7=3, 5"
Replace-Once $d.Content "@@PLACEHOLDER_25@@" "81÷9=9, 0"

Write-Host "Done"
